$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 4870
$ws1.Range("F11").Value = 1542
$ws1.Range("F15").Value = 1532
$ws1.Range("F20").Value = 4035
$ws1.Range("F21").Value = 4035
$ws1.Range("F23").Value = 3309
$ws1.Range("F24").Value = 767
$ws1.Range("F30").Value = 35

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 4870
$ws4.Range("F11").Value = 4870
$ws4.Range("F22").Value = 1542
$ws4.Range("F25").Value = 1532
$ws4.Range("F31").Value = 4035
$ws4.Range("F32").Value = 4035
$ws4.Range("F34").Value = 3309
$ws4.Range("F35").Value = 767
$ws4.Range("F41").Value = 35
